$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" ---
# This text shows up in the "Status" column of the per-locale sheets (zh-cn / de-de)
# and in the mirrored locale-status columns of the Overview sheet.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Narrow the "zh-cn"/"de-de" columns on Overview and the "Status" column on
#     the per-locale sheets from ~17.22 chars to ~13.41 chars ---

$wsOverview.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn)
$wsOverview.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de)

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
